{"js": "// Risk section: trim the placeholder/instructional text that trails the\n// \"T\u00ean r\u1ee7i ro:\" and \"M\u00f4 t\u1ea3 r\u1ee7i ro:\" labels, leaving just the labels\n// (work has been estimated, so the guidance text is no longer needed).\n\nconst body = context.document.body;\n\n// \"T\u00ean r\u1ee7i ro: ng\u1eafn g\u1ecdn trong 1 d\u00f2ng\" -> \"T\u00ean r\u1ee7i ro: \"\nconst nameHint = body.search(\"ng\u1eafn g\u1ecdn trong 1 d\u00f2ng\", { matchCase: true, matchWholeWord: false });\nnameHint.load(\"items\");\nawait context.sync();\nif (nameHint.items.length > 0) {\n  nameHint.items[0].getRange(Word.RangeLocation.whole).delete();\n  await context.sync();\n}\n\n// \"M\u00f4 t\u1ea3 r\u1ee7i ro: tool upload ch\u1eadm\" -> \"M\u00f4 t\u1ea3 r\u1ee7i ro: \"\nconst descHint = body.search(\" tool upload \", { matchCase: true, matchWholeWord: false });\ndescHint.load(\"items\");\nawait context.sync();\nif (descHint.items.length > 0) {\n  descHint.items[0].insertText(\" \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst slowWord = body.search(\"ch\u1eadm\", { matchCase: true, matchWholeWord: false });\nslowWord.load(\"items\");\nawait context.sync();\nif (slowWord.items.length > 0) {\n  slowWord.items[0].getRange(Word.RangeLocation.whole).delete();\n  await context.sync();\n}\n", "ps1": "# Risk section: trim the placeholder/instructional text that trails the\n# \"Ten rui ro:\" and \"Mo ta rui ro:\" labels, leaving just the labels\n# (work has been estimated, so the guidance text is no longer needed).\n\n$d = $word.ActiveDocument\n\n# \"Ten rui ro: ngan gon trong 1 dong\" -> \"Ten rui ro: \"\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"ng\u1eafn g\u1ecdn trong 1 d\u00f2ng\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found1) {\n    $rng1.Delete()\n}\n\n# \"Mo ta rui ro: tool upload cham\" -> \"Mo ta rui ro: \"\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\" tool upload \", $false, $false, $false, $false, $false, $true, 1, $false, \" \", 2)\nif ($found2) {\n    $rng2.Text = \" \"\n}\n\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"ch\u1eadm\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found3) {\n    $rng3.Delete()\n}\n"}
